# Move the "Senior Analyst - Myers Research" job block so it appears
# immediately before the "Research Director - PCCC" job block, and move
# the "Field Director - The Feldman Group" job block so it appears
# immediately before the "Programmer - Lake Research Partners" job block
# (i.e. swap the text content of each pair of blocks, since both blocks
# in each pair share the same paragraph shape: a Heading3 title, a
# subtitle paragraph, and three bullet paragraphs).

$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------
# Pair 1: "Research Director - PCCC" block  <->  "Senior Analyst - Myers
# Research" block
# ---------------------------------------------------------------------

Replace-Exact "Research Director - PCCC (Washington, DC) | August 2011 - August 2012" `
              "~~TMP1~~"
Replace-Exact "Senior Analyst - Myers Research (Austin, TX) | 2012 - 2014" `
              "Research Director - PCCC (Washington, DC) | August 2011 - August 2012"
Replace-Exact "~~TMP1~~" `
              "Senior Analyst - Myers Research (Austin, TX) | 2012 - 2014"

Replace-Exact "Political Research & Data Analysis (FLEEM System)" `
              "~~TMP2~~"
Replace-Exact "Political Research & Analysis" `
              "Political Research & Data Analysis (FLEEM System)"
Replace-Exact "~~TMP2~~" `
              "Political Research & Analysis"

Replace-Exact "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys" `
              "~~TMP3~~"
Replace-Exact "• Designed comprehensive survey instruments for specialized voting segments and niche markets" `
              "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys"
Replace-Exact "~~TMP3~~" `
              "• Designed comprehensive survey instruments for specialized voting segments and niche markets"

Replace-Exact "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren" `
              "~~TMP4~~"
Replace-Exact "• Developed sophisticated analytical products and reports that delivered actionable insights to clients" `
              "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"
Replace-Exact "~~TMP4~~" `
              "• Developed sophisticated analytical products and reports that delivered actionable insights to clients"

Replace-Exact "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver" `
              "~~TMP5~~"
Replace-Exact "• Co-developed a web application to manage all aspects of survey operations, from instrument design to data collection and analysis" `
              "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
Replace-Exact "~~TMP5~~" `
              "• Co-developed a web application to manage all aspects of survey operations, from instrument design to data collection and analysis"

# ---------------------------------------------------------------------
# Pair 2: "Programmer - Lake Research Partners" block  <->  "Field
# Director - The Feldman Group" block
# ---------------------------------------------------------------------

Replace-Exact "Programmer - Lake Research Partners (Austin, TX) | 2008" `
              "~~TMP6~~"
Replace-Exact "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012" `
              "Programmer - Lake Research Partners (Austin, TX) | 2008"
Replace-Exact "~~TMP6~~" `
              "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012"

Replace-Exact "Political Polling & Research" `
              "~~TMP7~~"
Replace-Exact "Political Campaign Management" `
              "Political Polling & Research"
Replace-Exact "~~TMP7~~" `
              "Political Campaign Management"

Replace-Exact "• Designed questionnaires and analyzed data for complex market research studies across diverse industries" `
              "~~TMP8~~"
Replace-Exact "• Managed all aspects of survey fielding for a multi-million dollar research firm, including scheduling, oversight, sampling, and quality control" `
              "• Designed questionnaires and analyzed data for complex market research studies across diverse industries"
Replace-Exact "~~TMP8~~" `
              "• Managed all aspects of survey fielding for a multi-million dollar research firm, including scheduling, oversight, sampling, and quality control"

Replace-Exact "• Conducted statistical modeling and analysis to address multifaceted consumer behavior questions" `
              "~~TMP9~~"
Replace-Exact "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings" `
              "• Conducted statistical modeling and analysis to address multifaceted consumer behavior questions"
Replace-Exact "~~TMP9~~" `
              "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"

Replace-Exact "• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps" `
              "~~TMP10~~"
Replace-Exact "• Created custom reports and data visualizations based on specific client requirements" `
              "• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"
Replace-Exact "~~TMP10~~" `
              "• Created custom reports and data visualizations based on specific client requirements"

Write-Host "Done swapping job blocks."
